$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the bold/centered/bordered style on B1 first
$cell1 = $ws.Range("B1")
$cell1.Font.Bold = $true
$cell1.HorizontalAlignment = -4108
$cell1.VerticalAlignment = -4160
$cell1.Borders.LineStyle = 1

# Copy the same formatting (without touching the value) to A2
$cell1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
